$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on all Price/Volume cells being updated so Excel
# does not silently reinterpret textual values (e.g. "1.00", "0.615") as numbers.
$textCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D46", "D49", "D50", "D51", "E2", "E3", "E4", "E5", "E6", "E7", "E8", "E9", "E10", "E11", "E12", "E13", "E14", "E15", "E16", "E17", "E18", "E19", "E20", "E21", "E22", "E23", "E24", "E25", "E26", "E28", "E29", "E30", "E31", "E32", "E33", "E34", "E35", "E36", "E37", "E39", "E40", "E41", "E42", "E43", "E44", "E45", "E46", "E47", "E48", "E49", "E50", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in the same order as the source diff.
$ws.Range("D2").Value = '69.598.07'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '3.544.33'
$ws.Range("E3").Value = '  -1.93%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '197.76'
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").Value = '586.05'
$ws.Range("E6").Value = '  -2.96%  '
$ws.Range("D7").Value = '0.615'
$ws.Range("E7").Value = '  -1.73%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.207'
$ws.Range("E9").Value = '  +0.82%  '
$ws.Range("E10").Value = '  -3.34%  '
$ws.Range("D11").Value = '52.14'
$ws.Range("E11").Value = '  -3.10%  '
$ws.Range("D12").Value = '0.0000290'
$ws.Range("E12").Value = '  -4.72%  '
$ws.Range("D13").Value = '9.35'
$ws.Range("E13").Value = '  -1.72%  '
$ws.Range("D14").Value = '680.59'
$ws.Range("E14").Value = '  +14.96%  '
$ws.Range("D15").Value = '4.103.28'
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("D16").Value = '69.648.43'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").Value = '12.46'
$ws.Range("E17").Value = '  -5.63%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.528.21'
$ws.Range("E18").Value = '  -2.08%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '18.63'
$ws.Range("E19").Value = '  -3.25%  '
$ws.Range("D20").Value = '0.122'
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("E21").Value = '  -2.39%  '
$ws.Range("D22").Value = '17.93'
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").Value = '108.17'
$ws.Range("E23").Value = '  +5.12%  '
$ws.Range("D24").Value = '5.26'
$ws.Range("E24").Value = '  +1.70%  '
$ws.Range("D25").Value = '4.41'
$ws.Range("E25").Value = '  -4.99%  '
$ws.Range("D26").Value = '2.95'
$ws.Range("E26").Value = '  -3.53%  '
$ws.Range("D28").Value = '10.35'
$ws.Range("E28").Value = '  -4.75%  '
$ws.Range("D29").Value = '9.71'
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("D30").Value = '33.49'
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("D31").Value = '4.38'
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("D32").Value = '6.93'
$ws.Range("E32").Value = '  -2.76%  '
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("D35").Value = '62.45'
$ws.Range("E35").Value = '  -1.20%  '
$ws.Range("D36").Value = '3.807.87'
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("D37").Value = '0.0₃0823'
$ws.Range("E37").Value = '  -5.91%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("D39").Value = '3.66'
$ws.Range("E39").Value = '  +3.23%  '
$ws.Range("D40").Value = '503.65'
$ws.Range("E40").Value = '  -3.66%  '
$ws.Range("D41").Value = '2.96'
$ws.Range("E41").Value = '  -7.34%  '
$ws.Range("D42").Value = '0.137'
$ws.Range("E42").Value = '  +2.14%  '
$ws.Range("E43").Value = '  -4.78%  '
$ws.Range("D44").Value = '35.01'
$ws.Range("E44").Value = '  -5.58%  '
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").Value = '2.96'
$ws.Range("E46").Value = '  +3.94%  '
$ws.Range("E47").Value = '  +1.67%  '
$ws.Range("E48").Value = '  -2.47%  '
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("D50").Value = '8.40'
$ws.Range("E50").Value = '  -2.51%  '
$ws.Range("D51").Value = '1.80'
$ws.Range("E51").Value = '  +20.95%  '
